$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the legend icon/label shared strings:
#   black square -> blue book, "noir" (black) -> "bleu" (blue)
#   orange book emoji swap
#   red book emoji swap
$ws.Cells.Replace("⬛", "📘", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("🟧", "📙", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("🟥", "📕", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("noir", "bleu", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
